$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (grandfather) and Row 5 (brother) Level column (C) change from "hard" to "easy"
$ws.Range("C4").Value = "easy"
$ws.Range("C5").Value = "easy"
